$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) "General/Program Chairs" rectangle: split " OAR Corporation, USA" so the
#        leading " OAR" becomes its own run reading "OAR " (Abhishek Dubey replaces
#        the old OAR Corp contact, the stray leading space is dropped in the process).
$chairs = $s.Shapes.Item(13)
$chairsRange = $chairs.TextFrame.TextRange
$oarPrefix = $chairsRange.Characters($chairsRange.Find(" OAR Corporation, USA").Start, 5)
$oarPrefix.Text = "OAR "
# This shape auto-fits its height to the text; re-apply the original height (553.7554pt
# = 7032694 EMU) so the rewrap caused by the edit above doesn't leave a stray size change.
$chairs.Height = 553.7555

# --- 2) Email textbox: fill in the real contact address instead of "TBA", and widen
#        the textbox so the longer email still fits (AutoFit growth captured in the diff).
$emailBox = $s.Shapes.Item(15)
$emailRange = $emailBox.TextFrame.TextRange
$tba = $emailRange.Find("TBA")
$tba.Text = "abhishek.dubey@Vanderbilt.Edu"

$emailBox.Width = 263.27625
